$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 12 & 13, plus two new header columns (D/E) for the rating_button
# choices. Shared-string entries must be created in this exact order so the
# sharedStrings.xml table comes out in the same sequence as upstream.
$ws.Range("A12").Value = "geopoint"
$ws.Range("B12").Value = "location"
$ws.Range("C12").Value = "Record your location"

$ws.Range("A13").Value = "rating_button"
$ws.Range("B13").Value = "agreement"
$ws.Range("C13").Value = "Do you agree that SurveySignal is cool & hip?"

$ws.Range("D1").Value = "choice1"
$ws.Range("E1").Value = "choice2"

$ws.Range("D13").Value = '<i class="fa fa-smile-o fa-2x"></i>'
$ws.Range("E13").Value = '<i class="fa fa-meh-o fa-2x"></i>'

$ws.Rows.Item(13).RowHeight = 45

[void]$ws.Range("E13").Select()
